{"js": "// Update the division-problem answers in the single table, cell by cell.\n// Each entry is addressed by its exact (row, column) position so that\n// duplicate/overlapping text values (e.g. \"42\u00f75=8, 2\" is both an old\n// value in one cell and a new value written into a different cell)\n// can never be confused with one another.\nconst updates = [\n  { row: 0, col: 0, oldText: \"67\u00f74=16, 3\", newText: \"72\u00f74=18, 0\" },\n  { row: 0, col: 1, oldText: \"60\u00f76=10, 0\", newText: \"65\u00f75=13, 0\" },\n  { row: 0, col: 2, oldText: \"19\u00f73=6, 1\", newText: \"35\u00f79=3, 8\" },\n  { row: 0, col: 3, oldText: \"73\u00f78=9, 1\", newText: \"83\u00f72=41, 1\" },\n  { row: 0, col: 4, oldText: \"50\u00f77=7, 1\", newText: \"33\u00f74=8, 1\" },\n\n  { row: 4, col: 0, oldText: \"27\u00f78=3, 3\", newText: \"30\u00f77=4, 2\" },\n  { row: 4, col: 1, oldText: \"24\u00f72=12, 0\", newText: \"37\u00f74=9, 1\" },\n  { row: 4, col: 2, oldText: \"84\u00f79=9, 3\", newText: \"59\u00f72=29, 1\" },\n  { row: 4, col: 3, oldText: \"57\u00f74=14, 1\", newText: \"42\u00f75=8, 2\" },\n  { row: 4, col: 4, oldText: \"30\u00f72=15, 0\", newText: \"84\u00f78=10, 4\" },\n\n  { row: 8, col: 0, oldText: \"74\u00f72=37, 0\", newText: \"72\u00f75=14, 2\" },\n  { row: 8, col: 1, oldText: \"42\u00f75=8, 2\", newText: \"30\u00f75=6, 0\" },\n  { row: 8, col: 2, oldText: \"21\u00f73=7, 0\", newText: \"38\u00f77=5, 3\" },\n  { row: 8, col: 3, oldText: \"61\u00f75=12, 1\", newText: \"58\u00f79=6, 4\" },\n  { row: 8, col: 4, oldText: \"14\u00f74=3, 2\", newText: \"49\u00f74=12, 1\" },\n\n  { row: 12, col: 0, oldText: \"10\u00f76=1, 4\", newText: \"20\u00f72=10, 0\" },\n  { row: 12, col: 1, oldText: \"35\u00f73=11, 2\", newText: \"60\u00f79=6, 6\" },\n  { row: 12, col: 2, oldText: \"83\u00f76=13, 5\", newText: \"99\u00f76=16, 3\" },\n  { row: 12, col: 3, oldText: \"49\u00f76=8, 1\", newText: \"73\u00f74=18, 1\" },\n  { row: 12, col: 4, oldText: \"28\u00f77=4, 0\", newText: \"45\u00f73=15, 0\" },\n\n  { row: 16, col: 0, oldText: \"65\u00f76=10, 5\", newText: \"50\u00f72=25, 0\" },\n  { row: 16, col: 1, oldText: \"49\u00f72=24, 1\", newText: \"54\u00f79=6, 0\" },\n  { row: 16, col: 2, oldText: \"69\u00f73=23, 0\", newText: \"76\u00f76=12, 4\" },\n  { row: 16, col: 3, oldText: \"53\u00f77=7, 4\", newText: \"76\u00f72=38, 0\" },\n  { row: 16, col: 4, oldText: \"87\u00f77=12, 3\", newText: \"92\u00f75=18, 2\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\n\nconst table = tables.items[0];\n\nfor (const { row, col, oldText, newText } of updates) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found in cell (${row}, ${col}): ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the division-problem answers in the single table, cell by cell.\n# Each entry is addressed by its exact 1-based (row, column) position so\n# that duplicate/overlapping text values (e.g. \"42\u00f75=8, 2\" is both an old\n# value in one cell and a new value written into a different cell) can\n# never be confused with one another, and a plain Find/Replace across the\n# whole document can't accidentally hit the wrong occurrence.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1;  Col = 1; Old = \"67\u00f74=16, 3\";  New = \"72\u00f74=18, 0\" },\n    @{ Row = 1;  Col = 2; Old = \"60\u00f76=10, 0\";  New = \"65\u00f75=13, 0\" },\n    @{ Row = 1;  Col = 3; Old = \"19\u00f73=6, 1\";   New = \"35\u00f79=3, 8\" },\n    @{ Row = 1;  Col = 4; Old = \"73\u00f78=9, 1\";   New = \"83\u00f72=41, 1\" },\n    @{ Row = 1;  Col = 5; Old = \"50\u00f77=7, 1\";   New = \"33\u00f74=8, 1\" },\n\n    @{ Row = 5;  Col = 1; Old = \"27\u00f78=3, 3\";   New = \"30\u00f77=4, 2\" },\n    @{ Row = 5;  Col = 2; Old = \"24\u00f72=12, 0\";  New = \"37\u00f74=9, 1\" },\n    @{ Row = 5;  Col = 3; Old = \"84\u00f79=9, 3\";   New = \"59\u00f72=29, 1\" },\n    @{ Row = 5;  Col = 4; Old = \"57\u00f74=14, 1\";  New = \"42\u00f75=8, 2\" },\n    @{ Row = 5;  Col = 5; Old = \"30\u00f72=15, 0\";  New = \"84\u00f78=10, 4\" },\n\n    @{ Row = 9;  Col = 1; Old = \"74\u00f72=37, 0\";  New = \"72\u00f75=14, 2\" },\n    @{ Row = 9;  Col = 2; Old = \"42\u00f75=8, 2\";   New = \"30\u00f75=6, 0\" },\n    @{ Row = 9;  Col = 3; Old = \"21\u00f73=7, 0\";   New = \"38\u00f77=5, 3\" },\n    @{ Row = 9;  Col = 4; Old = \"61\u00f75=12, 1\";  New = \"58\u00f79=6, 4\" },\n    @{ Row = 9;  Col = 5; Old = \"14\u00f74=3, 2\";   New = \"49\u00f74=12, 1\" },\n\n    @{ Row = 13; Col = 1; Old = \"10\u00f76=1, 4\";   New = \"20\u00f72=10, 0\" },\n    @{ Row = 13; Col = 2; Old = \"35\u00f73=11, 2\";  New = \"60\u00f79=6, 6\" },\n    @{ Row = 13; Col = 3; Old = \"83\u00f76=13, 5\";  New = \"99\u00f76=16, 3\" },\n    @{ Row = 13; Col = 4; Old = \"49\u00f76=8, 1\";   New = \"73\u00f74=18, 1\" },\n    @{ Row = 13; Col = 5; Old = \"28\u00f77=4, 0\";   New = \"45\u00f73=15, 0\" },\n\n    @{ Row = 17; Col = 1; Old = \"65\u00f76=10, 5\";  New = \"50\u00f72=25, 0\" },\n    @{ Row = 17; Col = 2; Old = \"49\u00f72=24, 1\";  New = \"54\u00f79=6, 0\" },\n    @{ Row = 17; Col = 3; Old = \"69\u00f73=23, 0\";  New = \"76\u00f76=12, 4\" },\n    @{ Row = 17; Col = 4; Old = \"53\u00f77=7, 4\";   New = \"76\u00f72=38, 0\" },\n    @{ Row = 17; Col = 5; Old = \"87\u00f77=12, 3\";  New = \"92\u00f75=18, 2\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $table.Cell($u.Row, $u.Col)\n    $cellRange = $cell.Range\n    $currentText = $cellRange.Text\n    if ($currentText -notmatch [regex]::Escape($u.Old)) {\n        throw \"Cell ($($u.Row), $($u.Col)) did not contain expected text '$($u.Old)' (found '$currentText')\"\n    }\n    $cellRange.Text = $u.New\n}\n"}
